$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '63.696.71'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.90%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.087.50'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -2.22%  '

# Row 4
$ws.Range('E4').Value = '  -0.11%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '565.66'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.71%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '159.48'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -4.88%  '

# Row 7
$ws.Range('E7').Value = '  +0.04%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.566'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -6.09%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '3.098.60'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.71%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.115'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.69%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.55'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -3.93%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.374'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -3.34%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.639.48'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -1.85%  '

# Row 14
$ws.Range('E14').Value = '  -2.27%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '63.845.86'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.83%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '24.33'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -3.63%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.102.25'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.91%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.0000152'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -2.44%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '401.60'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -3.47%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '5.15'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.42%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.21'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -4.75%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.89'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -3.44%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.02%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '66.99'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -3.57%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.474'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -4.69%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.191'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -5.98%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0₃0995'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -2.70%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.92'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.75%  '

# Row 29
$ws.Range('E29').Value = '  +0.31%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.78'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.08%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '20.95'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -3.44%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '163.07'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +4.96%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.80'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -4.60%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '6.14'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -3.28%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.10'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.22%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.33'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.91%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.63'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -3.20%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.563.97'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -4.80%  '

# Row 40
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '23.34'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.08%  '

# Row 41
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.05'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -3.57%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '37.82'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -3.13%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.678'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -5.45%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0607'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.40%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0251'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -3.77%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '5.10'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -6.47%  '

# Row 47
$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.999'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.22%  '

# Row 48
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '20.63'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.82%  '

# Row 49
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '278.54'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -3.35%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0963'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.50%  '

# Row 51
$ws.Range('E51').Value = '  +0.46%  '
